$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.452.71"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").Value = "3.543.22"
$ws.Range("E3").Value = "  -0.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.34"
$ws.Range("E5").Value = "  +6.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "555.37"
$ws.Range("E6").Value = "  -2.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("E7").Value = "  -2.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.664"
$ws.Range("E9").Value = "  -1.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.94"
$ws.Range("E10").Value = "  +11.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -3.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000274"
$ws.Range("E12").Value = "  +1.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.95"
$ws.Range("E13").Value = "  +1.10%  "

# Row 14
$ws.Range("D14").Value = "4.141.09"
$ws.Range("E14").Value = "  +0.15%  "

# Row 15
$ws.Range("D15").Value = "3.578.19"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16
$ws.Range("E16").Value = "  -0.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.79"
$ws.Range("E17").Value = "  +3.33%  "

# Row 18
$ws.Range("D18").Value = "67.698.93"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.99"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.05"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "392.29"
$ws.Range("E21").Value = "  -2.29%  "

# Row 22
$ws.Range("B22").Value = "RenderToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.21"
$ws.Range("E22").Value = "  +2.43%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  -1.80%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.37"
$ws.Range("E24").Value = "  -2.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  -2.80%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.81"
$ws.Range("E26").Value = "  +3.08%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  -2.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.97"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.93"
$ws.Range("E29").Value = "  -0.70%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "701.18"
$ws.Range("E30").Value = "  +9.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.37"
$ws.Range("E31").Value = "  -5.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.87"
$ws.Range("E32").Value = "  -2.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.44"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("E34").Value = "  -3.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.59"
$ws.Range("E35").Value = "  -3.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.413"
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("E38").Value = "  +0.84%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.137.22"
$ws.Range("E39").Value = "  -1.23%  "

# Row 40
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +25.89%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0714"
$ws.Range("E42").Value = "  -6.42%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.130"
$ws.Range("E43").Value = "  -2.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  -6.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.72"
$ws.Range("E45").Value = "  +7.64%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  -2.46%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.07"
$ws.Range("E47").Value = "  -1.58%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  -1.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.43"
$ws.Range("E49").Value = "  -1.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.67"
$ws.Range("E50").Value = "  -3.63%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  -1.55%  "
